# Scheduled market-data refresh: update price/profit columns (H-N) for the
# rows whose Universalis averages moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 8 - On the Drip
$ws.Range("H8").Value = 76.333336
$ws.Range("I8").Value = 76.333336
$ws.Range("K8").Value = 229.000008
$ws.Range("M8").Value = -90.00000800000001
# row 19 - Unbreak My Heart
$ws.Range("H19").Value = 3101
$ws.Range("I19").Value = 5074.1
$ws.Range("J19").Value = 1583.2307
$ws.Range("K19").Value = 5074.1
$ws.Range("L19").Value = 1583.2307
$ws.Range("M19").Value = -4899.1
$ws.Range("N19").Value = -1933.2307
# row 137 - Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 921.79486
$ws.Range("I137").Value = 881.8611
$ws.Range("J137").Value = 1401
$ws.Range("K137").Value = 2645.5833
$ws.Range("L137").Value = 4203
$ws.Range("M137").Value = -95.58329999999978
$ws.Range("N137").Value = -9303

$ws = $wb.Worksheets.Item("ARM")
# row 61 - Dealing with the Tough Stuff
$ws.Range("H61").Value = 2207.7778
$ws.Range("I61").Value = 1036.4286
$ws.Range("J61").Value = 3469.2307
$ws.Range("K61").Value = 1036.4286
$ws.Range("L61").Value = 3469.2307
$ws.Range("M61").Value = -824.4286
$ws.Range("N61").Value = -3893.2307
# row 63 - Rivets Run through It
$ws.Range("H63").Value = 3334933
$ws.Range("I63").Value = 5001249.5
$ws.Range("J63").Value = 2300
$ws.Range("K63").Value = 5001249.5
$ws.Range("L63").Value = 2300
$ws.Range("M63").Value = -5000563.5
$ws.Range("N63").Value = -3672
# row 66 - A Riveting Revival (L)
$ws.Range("H66").Value = 3334933
$ws.Range("I66").Value = 5001249.5
$ws.Range("J66").Value = 2300
$ws.Range("K66").Value = 25006247.5
$ws.Range("L66").Value = 11500
$ws.Range("M66").Value = -25002815.5
$ws.Range("N66").Value = -18364
# row 74 - As the Bolt Flies
$ws.Range("H74").Value = 1315.6
$ws.Range("I74").Value = 1550.1666
$ws.Range("K74").Value = 1550.1666
$ws.Range("M74").Value = -676.1666
# row 77 - Heavy Metal Banned (L)
$ws.Range("H77").Value = 1315.6
$ws.Range("I77").Value = 1550.1666
$ws.Range("K77").Value = 7750.833000000001
$ws.Range("M77").Value = -3382.833000000001
# row 97 - Ore for Me
$ws.Range("H97").Value = 360.78946
$ws.Range("I97").Value = 290.9375
$ws.Range("J97").Value = 733.3333
$ws.Range("K97").Value = 290.9375
$ws.Range("L97").Value = 733.3333
$ws.Range("M97").Value = 205.0625
$ws.Range("N97").Value = -1725.3333
# row 102 - Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1454.2858
$ws.Range("I102").Value = 1716
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 1716
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = -94
$ws.Range("N102").Value = -4044
# row 132 - Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3905.327
$ws.Range("I132").Value = 4292.3613
$ws.Range("J132").Value = 3034.5
$ws.Range("K132").Value = 12877.0839
$ws.Range("L132").Value = 9103.5
$ws.Range("M132").Value = -10347.0839
$ws.Range("N132").Value = -14163.5
# row 136 - Metal with Mettle
$ws.Range("H136").Value = 2207.7778
$ws.Range("I136").Value = 1036.4286
$ws.Range("J136").Value = 3469.2307
$ws.Range("K136").Value = 3109.2858
$ws.Range("L136").Value = 10407.6921
$ws.Range("M136").Value = -559.2857999999997
$ws.Range("N136").Value = -15507.6921

$ws = $wb.Worksheets.Item("BSM")
# row 134 - Ruthenium Supremium
$ws.Range("H134").Value = 41319.96
$ws.Range("I134").Value = 65363.688
$ws.Range("J134").Value = 2850
$ws.Range("K134").Value = 196091.064
$ws.Range("L134").Value = 8550
$ws.Range("M134").Value = -193556.064
$ws.Range("N134").Value = -13620

$ws = $wb.Worksheets.Item("CRP")
# row 31 - Wall Not Found
$ws.Range("H31").Value = 9526795
$ws.Range("I31").Value = 2323.9333
$ws.Range("J31").Value = 33337974
$ws.Range("K31").Value = 2323.9333
$ws.Range("L31").Value = 33337974
$ws.Range("M31").Value = -2028.9333
$ws.Range("N31").Value = -33338564
# row 34 - Armoires of the Rich and Famous
$ws.Range("H34").Value = 9526795
$ws.Range("I34").Value = 2323.9333
$ws.Range("J34").Value = 33337974
$ws.Range("K34").Value = 2323.9333
$ws.Range("L34").Value = 33337974
$ws.Range("M34").Value = -2121.9333
$ws.Range("N34").Value = -33338378
# row 58 - You Do the Heavy Lifting
$ws.Range("H58").Value = 1021.4737
$ws.Range("I58").Value = 965.17645
$ws.Range("K58").Value = 965.17645
$ws.Range("M58").Value = -762.17645
# row 132 - Hull Lotta Damage
$ws.Range("H132").Value = 1799.591
$ws.Range("I132").Value = 1044
$ws.Range("K132").Value = 3132
$ws.Range("M132").Value = -602
# row 134 - Wood You Be Quiet
$ws.Range("H134").Value = 1153.32
$ws.Range("I134").Value = 1060.591
$ws.Range("J134").Value = 1833.3334
$ws.Range("K134").Value = 3181.773
$ws.Range("L134").Value = 5500.0002
$ws.Range("M134").Value = -646.7729999999997
$ws.Range("N134").Value = -10570.0002
# row 136 - Turali Quality
$ws.Range("H136").Value = 1021.4737
$ws.Range("I136").Value = 965.17645
$ws.Range("K136").Value = 2895.52935
$ws.Range("M136").Value = -345.5293500000002

$ws = $wb.Worksheets.Item("GSM")
# row 19 - Better Four Eyes than None
$ws.Range("H19").Value = 10001525
$ws.Range("I19").Value = 40000000
$ws.Range("J19").Value = 2033.3334
$ws.Range("K19").Value = 40000000
$ws.Range("L19").Value = 2033.3334
$ws.Range("M19").Value = -39999712
$ws.Range("N19").Value = -2609.3334
# row 21 - Forever 21K
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 35000
$ws.Range("L21").Value = 35000
$ws.Range("N21").Value = -35346
# row 30 - Dog Tags Are for Dogs
$ws.Range("H30").Value = 35000
$ws.Range("J30").Value = 35000
$ws.Range("L30").Value = 35000
$ws.Range("N30").Value = -35210
# row 80 - Needs More Prayerbell
$ws.Range("H80").Value = 4584.95
$ws.Range("I80").Value = 3777.7778
$ws.Range("J80").Value = 5245.364
$ws.Range("K80").Value = 3777.7778
$ws.Range("L80").Value = 5245.364
$ws.Range("M80").Value = -2779.7778
$ws.Range("N80").Value = -7241.364
# row 83 - With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 4584.95
$ws.Range("I83").Value = 3777.7778
$ws.Range("J83").Value = 5245.364
$ws.Range("K83").Value = 18888.889
$ws.Range("L83").Value = 26226.82
$ws.Range("M83").Value = -13896.889
$ws.Range("N83").Value = -36210.82
# row 132 - On Board for Lar
$ws.Range("H132").Value = 51420.27
$ws.Range("I132").Value = 57900.918
$ws.Range("J132").Value = 4759.6
$ws.Range("K132").Value = 173702.754
$ws.Range("L132").Value = 14278.8
$ws.Range("M132").Value = -171172.754
$ws.Range("N132").Value = -19338.8

$ws = $wb.Worksheets.Item("LTW")
# row 11 - A Thorn in One's Hide
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = ""
# row 46 - Supply Side Logic
$ws.Range("H46").Value = 1032.1177
$ws.Range("I46").Value = 266.66666
$ws.Range("J46").Value = 1449.6364
$ws.Range("K46").Value = 266.66666
$ws.Range("L46").Value = 1449.6364
$ws.Range("M46").Value = -78.66665999999998
$ws.Range("N46").Value = -1825.6364
# row 61 - Spelling Me Softly
$ws.Range("H61").Value = 1366.4445
$ws.Range("I61").Value = 1328.2858
$ws.Range("K61").Value = 1328.2858
$ws.Range("M61").Value = -1126.2858
# row 113 - Peace in Rest
$ws.Range("H113").Value = 1366.4445
$ws.Range("I113").Value = 1328.2858
$ws.Range("K113").Value = 1328.2858
$ws.Range("M113").Value = 841.7141999999999
# row 132 - Tenets of Tanning
$ws.Range("H132").Value = 2340.7727
$ws.Range("I132").Value = 2407.6155
$ws.Range("K132").Value = 7222.8465
$ws.Range("M132").Value = -4692.8465
# row 136 - Respect for Br'aax
$ws.Range("H136").Value = 7577.8887
$ws.Range("I136").Value = 10609.272
$ws.Range("K136").Value = 31827.816
$ws.Range("M136").Value = -29277.816

$ws = $wb.Worksheets.Item("WVR")
# row 13 - Time for Acton
$ws.Range("H13").Value = 527.5
$ws.Range("I13").Value = 527.5
$ws.Range("K13").Value = 527.5
$ws.Range("M13").Value = -387.5
# row 92 - Modest Beginnings
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
# row 107 - Flax Wax
$ws.Range("H107").Value = 213.33333
$ws.Range("I107").Value = 223
$ws.Range("J107").Value = 194
$ws.Range("K107").Value = 669
$ws.Range("L107").Value = 582
$ws.Range("M107").Value = 1251
$ws.Range("N107").Value = -4422
# row 132 - Comfy Cabins
$ws.Range("H132").Value = 2034.525
$ws.Range("I132").Value = 1719.72
$ws.Range("J132").Value = 2559.2
$ws.Range("K132").Value = 5159.16
$ws.Range("L132").Value = 7677.599999999999
$ws.Range("M132").Value = -2629.16
$ws.Range("N132").Value = -12737.6
# row 136 - Weaving the Envelope
$ws.Range("H136").Value = 4700.485
$ws.Range("I136").Value = 6111.0835
$ws.Range("J136").Value = 938.8889
$ws.Range("K136").Value = 18333.2505
$ws.Range("L136").Value = 2816.6667
$ws.Range("M136").Value = -15783.2505
$ws.Range("N136").Value = -7916.6667
